# Update the "Förändrad" (Changed) date column (C) for every data row.
# The value 45203 (2023-10-04) is bumped by one day to 45204 (2023-10-05)
# for all rows that currently hold it, mirroring an automated "last
# updated" timestamp refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
